# Inserts two new data rows (new rows 90 and 91) into the Ají price table,
# pushing the existing rows 90-163 down to 92-165 and extending the sheet's
# used range from A1:R163 to A1:R165.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above the current row 90 (shifts 90:163 -> 92:165).
$ws.Rows("90:91").Insert()

# --- New row 90: Ají / Americana (o) / Primera, Provincia de Limarí ---
$ws.Range("A90").Value = 8
$ws.Range("B90").Value = "Terminal La Palmera de La Serena"
$ws.Range("C90").Value = "Coquimbo"
$ws.Range("D90").Value = 44566
$ws.Range("E90").Value = 4
$ws.Range("F90").Value = 100112021
$ws.Range("G90").Value = "Ají"
$ws.Range("H90").Value = "Americana (o)"
$ws.Range("I90").Value = "Primera"
$ws.Range("J90").Value = 500
$ws.Range("K90").Value = 26000
$ws.Range("L90").Value = 27000
$ws.Range("M90").Value = 26500
$ws.Range("N90").Value = "$/caja 25 kilos"
$ws.Range("O90").Value = "Provincia de Limarí"
$ws.Range("P90").Value = 1060
$ws.Range("Q90").Value = 25
$ws.Range("R90").Value = "Hortaliza"

# --- New row 91: Ají / Inferno / Primera, Provincia de Limarí ---
$ws.Range("A91").Value = 8
$ws.Range("B91").Value = "Terminal La Palmera de La Serena"
$ws.Range("C91").Value = "Coquimbo"
$ws.Range("D91").Value = 44566
$ws.Range("E91").Value = 4
$ws.Range("F91").Value = 100112021
$ws.Range("G91").Value = "Ají"
$ws.Range("H91").Value = "Inferno"
$ws.Range("I91").Value = "Primera"
$ws.Range("J91").Value = 600
$ws.Range("K91").Value = 16000
$ws.Range("L91").Value = 17000
$ws.Range("M91").Value = 16500
$ws.Range("N91").Value = "$/caja 15 kilos"
$ws.Range("O91").Value = "Provincia de Limarí"
$ws.Range("P91").Value = 1100
$ws.Range("Q91").Value = 15
$ws.Range("R91").Value = "Hortaliza"
